$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = '5 Şubattan beri yayınlanan ders programından dün ders çıkarıldı ve programlarımız onaylanmışken oldu bu. Yeniden sistem açtırıp ders seçmek zorunda kaldık. Bu durumda birçok derste kota sorunu ve çakışma ile karşılaştık. Diyecek en ufak bir sözüm bile yok sana ışık üniversitesi..'
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = 'Gecen donem online diye acilmayan lab derslerini bu donem oyle bi ayarlamislari ki her dersle cakisiyor ve okul 1 dersten fazla cakismaya izin vermedigi icin dersleri alamiyoruz okulumuz 1 sene daha uzuyor.Herkes magdur bu durumdan dersleri acilmamasi okulun sorunu degilmis gibi birde destekte olmuyorlar ogrenciye.Bu sekilde magdur olan arkadaslarimiz varsa gonderiyi begensin topluca rektorluge mail atalim yeter artik bu okul nasil biticek yoksa?!'
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = 'Bir dönem boyunca her dersine girdiğim, hiçbir guizini kaçırmadığım ve bütün ödevleri eksiksiz tamamladığım dersten sırf sunum yapmadığım için kaldığımı öğrendim. Dönem boyunca harcadığım zamanımı ve emeğimi çöpe attılar resmen. Işık üniversitesi allah belanı versin gençliğimizi çürüttün.'
$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(33, 2).Value = 'Her gecen gun bilgisayar muhendisligi hocalari daha da dusuncesiz davranamazlar diyorum ve her seferinde yaniliyorum. SOFT dersi icin gonderilen sinav yonetmeligini bir gorseniz yani bu baski altinda bilende unutur zaten yapamaz. "Sozde" calisanin hakkini korumak icin alinan onlemler gercegi yansitmiyor bu yuzden umarim en kisa zamanda okullar acilir ve bildigimiz derslerden kalmayiz.'
$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 2).Value = 'Okulumuzda öyle bir eğitim anlayisi var ki hocanın şakşakcısı olursan tüm dersleri en kötü ba ile geçerken, gruptan aykırı bir şeyler başaran biri olursan dc dd gibi notlarla geçiyorsun. Sorarım sizlere biz şak şak için mi eğitim alıyoruz, başarı için mi?'
$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(35, 2).Value = 'Sınavı respondusta yaparken 15 dakika içinde çöktü bilgisayar. Daha ne yapmamız gerekiyo derdimizi anlamaları için.'
$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 2).Value = 'Merhaba herkese, Tıpkı maillerdeki gibi başlamak istedim. Görünüşte ilgili ve bir o kadarda fikirlerimize duyarlı bir okul gibi gözükse de gerçek sadece mağduriyetten ibaret ENGL1301 sınavı tam bir rezaletle ile karşı karşıyaydık bu okul bana ilk dönemim olmasına rağmen stres sıkıntıdan başka birşey vermedi Işık üniversitesi hayatımı aydınlatmadan ziyade hayatımı kararttı. Haksızlık bu respondusss yerine çözüm bulunmalı. 2 dönem umarım yüzyüze egitim olur.'
$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 2).Value = 'Makine mühendisliği okuyup statik dersinin finaline giren arkadaşlarla aynı problemi yaşadık klasik soru sorup zoomdan kamera açtırıp her soru için 2 3 şıklı yapıp ve her soru için 10 dk süre verip fotoğrafını çekip sisteme yükleme gibi bi eziyet inanınki insanlığa sığmaz ve bu dersi verememek okulu 1 sene uzatmak demek. Şikayette bulunsak derse giren hocaya ters düştüğümüz taktirde okul olcak 8 sene başka üniversitede okuyan arkadaşlarım böyle bi sınav çeşidi olmaz derken biz nie hala bu eziyetlere katlanmak zorunda kalıyoruz anlamıyorum yazıklar olsun diyip geçicem.'
$ws.Cells.Item(38, 1).Value = 36
$ws.Cells.Item(38, 2).Value = 'Respondus kullanılan ekonomi sınavında 45 soruya 60 dakika vermek. Kolumuzu kıpırdatamıyoruz zaten neden bu kadar soru sayısı anlamış değilim ?'
$ws.Cells.Item(39, 1).Value = 37
$ws.Cells.Item(39, 2).Value = 'Bir üniversite düşün ki her hoca sadece kendi dersi varmışçasına ödevler projeler veriyor. Bir dersin bir dönemde 2 proje ,4 guiz, 2 vize, 1 finali olursa nasıl yetişecek diğer dersler? Birisi de ödevi el yazısıyla istiyor hani word kullanmak kolay ya.Hayret ediyorum. Amatörlük mü,art niyet mi karar sizin.'
$ws.Cells.Item(40, 1).Value = 38
$ws.Cells.Item(40, 2).Value = 'Okulumuzdaki bazı hocalar,bu bazı hocalara şanslıyım ki benim dersine girdiklerim dahil değil ama bu buraya bir şey yazmamam için bi sebep değil.Özellikle malum fakültedeki öğrencilere ne yapsak da ızdırap çektirsek diye düşündüklerini düşüneceğiz yakında.Ya her şeyden önce insan bu öğrenciler.Tamam zorluyorsunuz bu bir yandan iyi bir şey onlar için ama bir yandan da insanların ailesinde hasta olanlar uğraşanlar edenler var,donanım problemi yaşayanlar var,malum programın bilgisayarı aptal etmesinden dolayı sorun yaşayanlar var ve bu uygulamayı guizler sırasında istiyorlar guiz esnasında dondu kaldı tanıdıgımın bilgisayarı.Neden bu eziyet neden? Zaten bir sürü para ödüyoruz artık gerçekten insanları derste bırakmaya uğraşıyorlar! Ben dediğm gibi öğrencilerine bu tür zulümler uygulayan desek yanlış olmaz onlarla aynı derslere denk gelmiyorum AMA bir sürü öğrenciye yazık günah.Gerçekten hoşgörülü olan canım hocalarıma da teşekkürlerimi sunuyorum.Zor dönemden geçen bir insanız hepimiz.'
$ws.Cells.Item(41, 1).Value = 39
$ws.Cells.Item(41, 2).Value = 'Sınav kağıdımı başka bir üniversitede hoca olan bir tanıdığıma gönderdim ve notlandırmasını rica ettim. İşlem hatasına göz yumamadığını ve 68 verdiğini söyledi. Bizim okuldaki hoca 19 vermiş 19. Bazı hocalar cidden kafayı yemiş durumda. Aile büyüklerimin hepsi korona nerdeyse. Hepsine ben yetişmeye çalışıyorum. Babanem düşüp omuzunu kırdı. Hareket edemiyor herşeyine ben koşuyorum. Annanemle dedem korona ve ağır geçiriyorlar. Dedem yoğun bakımda. Annanemede ben bakıyorum. Alzheimer hastası. Her saat dedemin nerde olduğu halkında yalan söyleyip tansiyonunu düşürmeye sakinleştirmeye çalışıyorum. Maddi durumumuz da yeterli değil. Kendim internetten bi şekilde para kazanmaya çalışıyorum ihtiyaçlarımızı karşılamak için. Bunların hiç birinden şikayetçi değilim. İsyan da etmiyorum. İsyan ettiğim şey bukadar karmaşanın içinde ders çalışmaya vakit ayırmaya çalışıyorum. Çünkü seneye bursum kesilecek ve ödeyemem.'
$ws.Cells.Item(42, 1).Value = 40
$ws.Cells.Item(42, 2).Value = 'Zaten çok zor bir süreç yaşarken online eğitimle cebelleşirken neredeyse bütün akademisyenlerin Dersleri anlatmaya bile yeltenmeyip Google''dan bakın internetten bakın araştırın ve öyle çizin öğrenin demesi artık kafayı yememe sebep oluyor biz öğrenciler olarak bu süreçte iyisiyle kötüsüyle belki de en çok derslere fokuslandık ama akademisyenler her şey normalmiş gibi itirazlarımızı bile tembelliğe vuruyor ayrıca hiçbir şekilde gerçekten öğrenmek istediğimiz için ve anlatmaya tenezzül etmedikleri için sitem ettiğimizi fark edemiyorlar Artık bıktım tüm sorunlarla boğuşurken derslerden kalma korkusuyla öğrenmeden kopyalayarak çizim yapmaktan da bıktım beni buna zorlayan sistemden de bıktım lütfen bazı şeylere tepki göstermekten çekinmeyin 1-2 kişi tepki gösterince bunu tembel öğrencilerin ses çıkarması olarak algılıyorlar ama çoğunluk tepki gösterirse gerçekten değiştiremeyeceğimiz hiçbir şey yok ışık öğrencisinin bunun farkına varması gerekiyor birlikten kuvvet doğar.'
$ws.Cells.Item(43, 1).Value = 41
$ws.Cells.Item(43, 2).Value = 'Merhaba kolay gelsin, şu konuya kimse değinmiyor ama hocaların çoğunluğu derse haber vermeden geç geliyorlar. Yeri gelince örgünden farkı yok artık alıştınız diyen hocalarımızdan ricam bir zahmet örgünde ki gibi vaktinde yataktan çıkıp derse gelmeleri. Gerçekten PC başında saatler harcarken artık 1 dakikasının bile bir hocayı beklemekle geçmesine katlanamaz oldum.'
$ws.Cells.Item(44, 1).Value = 42
$ws.Cells.Item(44, 2).Value = 'Bu okulda hakkınızı aramaya çalışmayın arkadaşlar boşa çabalarsınız hiç bişey değiştirmez çabalarınız sınavları geç bitirme projelerinizde bile o kadar büyük hakkınızı yiyorlar ki aklınız almıyo kişiye göre ayrımcılık desen o da var işleri yarım insanlar yüksek not alırken tam anlamıyla herşeyini tamamlamış insanların en düşük notu alması gibi hakkını savunmaya çalıştıgında dekanlığa bile gitsen yapılanları anlatsan hakkını arasan hiç bişey değişmiyo ne de olsa kendi hocaları onları savunacaklar öğrenciyi kim savunsun öğrenci dersten kalsın okula paralar gelsin onların en çok istiyeceği şey o yüzden çabalamayınn arkadaşlarr daha çok üzüleceksiniz bu okulda.'
$ws.Cells.Item(45, 1).Value = 43
$ws.Cells.Item(45, 2).Value = 'Dersin ismini vermeyeceğim. İlk defa başıma geliyor ve yazmak istedim. Bu sıradan bir serzeniş olmayacak. Yazının devamını yazdım fakat sildim. Mesele sınav kuralları ile alakalıydı. Her şeyi geçtim, sadece şu konuda çok üzgünüm. Biz öğrenciler hakkımızı arayamıyoruz herhangi bir konuda. Okul dediğimiz bina bize teoriden başka hiçbir şey öğretmiyor. Bize hayata dair bir şeyler vermeli. Hakkımızı aramamız gerektiğini hatırlatmalı fakat başımıza bir şey gelir diye susuyoruz. Okul bizi, hayata karşı daha korkak hale getiriyor ama farkında değiliz. Bu sadece bu okul nezdinde değil. Türkiye''deki bütün okullar böyle. Benim çok başıma geldi. Okuldayken hocanın karşısında iki büklüm oldum defalarca. Bu korku düzeni değişmeli.'
$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 2).Value = 'Derste respondus programının bilgisayara zarar verdiğini, bazı bilgisayarların çökmesine neden olduğunu söyledim ve çok komik bir yanıt aldım hocamız abartmayın madem çöküyorsa verilerinizi başka bir yere yedekleyin sorun ortadan kalkmış olur dedi, herkesin durumu olmayabiliyor benim için verilerden çok bilgisayarım daha önemli bilgisayarımda bir sorun oluştuğunda sanki aynı gün tekrar bilgisayar alabilecekmişim gibi davranıldı. Hocaların kendileri, öğrencilerden daha çok teknik problem yaşarken bizlerden bu kadar şeyi aynı anda yerine getirilmesini zorunlu tutmaları çok anlayışsızca eğitimci olarak ilk başta anlayışlı olmaları ve empati kurabilmelerini diliyorum.'
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 2).Value = 'Merhaba respondus artık gerçekten çalışmamaya başladı, bu zamana kadar takır takır çalışıyordu fakat 3 gündür uygulamada sorun yaşıyorum. Respondus uygulaması beni sınava almadı ilk 8 dakika, bütün kamera ve mikrofonu tanımladı yani bütün aşamalara onay verdi fakat en sonda sınavı başlatmadı, yüzümü tanımlamadığını belirtti her açıdan denememe rağmen birkaç dakika bu nedenden ötürü respondusu kapatıp tekrar yaptım. Bu yakın bir zamanda sizin de başınıza gelebilir. Lütfen artık birlik olup bu durumu düzeltelim.'
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 2).Value = 'Bir şey söylemek istiyorum, yaz okulu dönemi respondus ve webcam eklentisi hakkında ben arkadaşlarımıza uzuncana açıklama yaparken bize "çalış dersi geç", "ağlama", "okul okuyorsan yapacaksın" vb tepkiler veren arkadaşların bir kısmı dahil olmak üzere, bizi sallamayı geçtim bir de dalga geçenler dahil tutuşmuş şekilde çözüm arıyorlar, gerçekten gülünç durumdasınız. He belirteyim, aldığım derslerde respondus isteyen yok şu an ancak gelin şu durum beter hal almadan düzgün bir çözüm bulalım, hep birlilte. Tweet atmakla ses duyuramadık, cimer de çok ciddiye almıyor açıkçası, ya çözeceğiz ya da geçmiş olsun arkadaşlar.'
$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 2).Value = 'Böyle bir saçmalık yok arkadaşlar hepimiz bir olalım ve şu kamera sistemini lütfen iptal ettirelim bahsettiğim sınav İngilizce sınavı ve hocalar acımasızca bize kameranı yoksa kamera alın diyor bu nasıl bir tavır belki benim durumum yok nerden biliyorsun madem öyle sevgili hocalarımız bize kamera alsın ben özel hayatıma müdahale edilmesini istemiyorum sonuç itibari ile uzaktan eğitim bizim isteğimizle olmadı bu kadar sıkıcaklarına bizi okula çağırıp gerekli önlemleri alıp bizleri okulda sınavda yapsınlar arkadaşlar gelin twitinde ıstegramda birlik olup şu işlemi iptal ediyorum ayrıca ben çalışıyorum 1:30 saat sınavmı olur nasıl yapacağımı kara kara düşünüyorum lütfen bir olalım bize destek çıkın lütfen bugün bize yarın size büyütelim bu olayı.'
$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(50, 2).Value = 'Merhaba respondus uygulaması ile ilgili bir şey belirtmek istiyorum. Sınav saati normalde 30 dakika, fakat respondusa girmeye çalıştığım sırada beni sınava almadı. "Yoğunluktan dolayı sizi sınava alamıyoruz, lütfen birkaç dakika bekleyin" yazısı ile karşı karşıya kaldım. Sınavı daha kısa sürede yapmak zorunda kaldım. Ayrıca bilgisayara yeni format atmama rağmen respondusu yükledikten sonra bilgisayarım yavaşladı ve de bazı uygulamalar çalışmamaya başladı. Lütfen bu konu hakkında birlik olup cimere yazalım.'

$ws.Range("G53").Select()
$excel.ActiveWindow.ScrollRow = 33

Write-Output "Done"